# Applies the price/volume updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.590.24"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "1.754.17"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").Value = "  +0.06%  "

$origStyleD5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.25"
$ws.Range("D5").Style = $origStyleD5
$ws.Range("E5").Value = "  -0.14%  "

$ws.Range("E6").Value = "  +0.03%  "

$origStyleD7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4499"
$ws.Range("D7").Style = $origStyleD7
$ws.Range("E7").Value = "  +4.42%  "

$origStyleD8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3551"
$ws.Range("D8").Style = $origStyleD8
$ws.Range("E8").Value = "  -1.52%  "

$origStyleD9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07450"
$ws.Range("D9").Style = $origStyleD9
$ws.Range("E9").Value = "  -1.80%  "

$origStyleD10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.61"
$ws.Range("D10").Style = $origStyleD10
$ws.Range("E10").Value = "  -1.44%  "

$ws.Range("E11").Value = "  -2.72%  "

$ws.Range("E12").Value = "  +0.03%  "

$origStyleD13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.71"
$ws.Range("D13").Style = $origStyleD13
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("E14").Value = "  -1.56%  "

$ws.Range("E15").Value = "  -1.27%  "

$ws.Range("D16").Value = "1.752.90"
$ws.Range("E16").Value = "  -0.41%  "

$origStyleD17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.42"
$ws.Range("D17").Style = $origStyleD17
$ws.Range("E17").Value = "  +1.15%  "

$origStyleD18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001056"
$ws.Range("D18").Style = $origStyleD18
$ws.Range("E18").Value = "  -1.12%  "

$origStyleD19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06456"
$ws.Range("D19").Style = $origStyleD19
$ws.Range("E19").Value = "  +0.30%  "

$origStyleD21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.06"
$ws.Range("D21").Style = $origStyleD21
$ws.Range("E21").Value = "  -0.26%  "

$origStyleD22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.740"
$ws.Range("D22").Style = $origStyleD22
$ws.Range("E22").Value = "  -2.40%  "

$ws.Range("D23").Value = "27.637.69"
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("E24").Value = "  -0.59%  "

$origStyleD25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.091"
$ws.Range("D25").Style = $origStyleD25
$ws.Range("E25").Value = "  +0.28%  "

$origStyleD26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.61"
$ws.Range("D26").Style = $origStyleD26
$ws.Range("E26").Value = "  +1.40%  "

$origStyleD27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.17"
$ws.Range("D27").Style = $origStyleD27
$ws.Range("E27").Value = "  -1.79%  "

$ws.Range("D28").Value = "1.952.56"
$ws.Range("E28").Value = "  -0.39%  "

$origStyleD29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.089"
$ws.Range("D29").Style = $origStyleD29
$ws.Range("E29").Value = "  -3.11%  "

$origStyleD30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.42"
$ws.Range("D30").Style = $origStyleD30
$ws.Range("E30").Value = "  -0.26%  "

$origStyleD31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.099"
$ws.Range("D31").Style = $origStyleD31
$ws.Range("E31").Value = "  -0.08%  "

$origStyleD32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09197"
$ws.Range("D32").Style = $origStyleD32
$ws.Range("E32").Value = "  +2.76%  "

$origStyleD33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.663"
$ws.Range("D33").Style = $origStyleD33
$ws.Range("E33").Value = "  -0.76%  "

$origStyleD34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.495"
$ws.Range("D34").Style = $origStyleD34
$ws.Range("E34").Value = "  -1.97%  "

$ws.Range("E35").Value = "  -0.72%  "

$origStyleD36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.70"
$ws.Range("D36").Style = $origStyleD36
$ws.Range("E36").Value = "  -4.38%  "

$origStyleD37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06031"
$ws.Range("D37").Style = $origStyleD37
$ws.Range("E37").Value = "  +0.21%  "

$origStyleD38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2083"
$ws.Range("D38").Style = $origStyleD38
$ws.Range("E38").Value = "  -1.60%  "

$origStyleD39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6287"
$ws.Range("D39").Style = $origStyleD39
$ws.Range("E39").Value = "  -1.22%  "

$origStyleD40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.927"
$ws.Range("D40").Style = $origStyleD40
$ws.Range("E40").Value = "  -0.56%  "

$origStyleD41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.182"
$ws.Range("D41").Style = $origStyleD41
$ws.Range("E41").Value = "  -0.42%  "

$ws.Range("E42").Value = "  -0.59%  "

$origStyleD43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.762"
$ws.Range("D43").Style = $origStyleD43
$ws.Range("E43").Value = "  -1.92%  "

$origStyleD44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.19"
$ws.Range("D44").Style = $origStyleD44
$ws.Range("E44").Value = "  -0.88%  "

$ws.Range("E45").Value = "  +0.08%  "

$origStyleD46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5855"
$ws.Range("D46").Style = $origStyleD46
$ws.Range("E46").Value = "  -1.31%  "

$origStyleD47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.33"
$ws.Range("D47").Style = $origStyleD47
$ws.Range("E47").Value = "  -0.31%  "

$ws.Range("E48").Value = "  -2.74%  "

$origStyleD49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06886"
$ws.Range("D49").Style = $origStyleD49
$ws.Range("E49").Value = "  +0.17%  "

$origStyleD50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.127"
$ws.Range("D50").Style = $origStyleD50
$ws.Range("E50").Value = "  -3.67%  "

$origStyleD51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.62"
$ws.Range("D51").Style = $origStyleD51
$ws.Range("E51").Value = "  -2.16%  "
